$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '28.029.21'
$ws.Range("E2").Value = '  -0.23%  '
$ws.Range("D3").Value = '1.815.38'
$ws.Range("E3").Value = '  +2.25%  '
$ws.Range("E4").Value = '  -0.12%  '
$ws.Range("D5").NumberFormat = '@'
$ws.Range("D5").Value = '336.70'
$ws.Range("E5").Value = '  -0.67%  '
$ws.Range("D6").NumberFormat = '@'
$ws.Range("D6").Value = '0.9997'
$ws.Range("E6").Value = '  -0.23%  '
$ws.Range("D7").NumberFormat = '@'
$ws.Range("D7").Value = '0.4268'
$ws.Range("E7").Value = '  +11.76%  '
$ws.Range("D8").NumberFormat = '@'
$ws.Range("D8").Value = '0.3513'
$ws.Range("E8").Value = '  +2.89%  '
$ws.Range("D9").NumberFormat = '@'
$ws.Range("D9").Value = '45.71'
$ws.Range("E9").Value = '  -2.69%  '
$ws.Range("E10").Value = '  +0.43%  '
$ws.Range("D11").NumberFormat = '@'
$ws.Range("D11").Value = '0.07444'
$ws.Range("E11").Value = '  +0.94%  '
$ws.Range("D12").NumberFormat = '@'
$ws.Range("D12").Value = '22.91'
$ws.Range("E12").Value = '  -0.81%  '
$ws.Range("E13").Value = '  +0.04%  '
$ws.Range("D14").NumberFormat = '@'
$ws.Range("D14").Value = '6.258'
$ws.Range("E14").Value = '  -1.85%  '
$ws.Range("D15").NumberFormat = '@'
$ws.Range("D15").Value = '7.304'
$ws.Range("E15").Value = '  -1.15%  '
$ws.Range("D16").Value = '1.813.10'
$ws.Range("E16").Value = '  +2.01%  '
$ws.Range("D17").NumberFormat = '@'
$ws.Range("D17").Value = '0.00001085'
$ws.Range("E17").Value = '  +0.82%  '
$ws.Range("D18").NumberFormat = '@'
$ws.Range("D18").Value = '0.06688'
$ws.Range("E18").Value = '  +0.50%  '
$ws.Range("D19").NumberFormat = '@'
$ws.Range("D19").Value = '82.02'
$ws.Range("E19").Value = '  -0.58%  '
$ws.Range("D20").NumberFormat = '@'
$ws.Range("D20").Value = '1.000'
$ws.Range("E20").Value = '  -0.09%  '
$ws.Range("E21").Value = '  -0.63%  '
$ws.Range("D22").NumberFormat = '@'
$ws.Range("D22").Value = '6.439'
$ws.Range("E22").Value = '  +0.83%  '
$ws.Range("D23").Value = '28.064.82'
$ws.Range("E23").Value = '  -0.15%  '
$ws.Range("D24").NumberFormat = '@'
$ws.Range("D24").Value = '11.87'
$ws.Range("E24").Value = '  -1.47%  '
$ws.Range("D25").NumberFormat = '@'
$ws.Range("D25").Value = '2.374'
$ws.Range("E25").Value = '  -0.19%  '
$ws.Range("D26").NumberFormat = '@'
$ws.Range("D26").Value = '2.488'
$ws.Range("E26").Value = '  +3.48%  '
$ws.Range("D27").NumberFormat = '@'
$ws.Range("D27").Value = '20.70'
$ws.Range("E27").Value = '  +0.05%  '
$ws.Range("D28").NumberFormat = '@'
$ws.Range("D28").Value = '155.42'
$ws.Range("E28").Value = '  +1.32%  '
$ws.Range("D29").Value = '2.020.54'
$ws.Range("E29").Value = '  +2.06%  '
$ws.Range("D30").NumberFormat = '@'
$ws.Range("D30").Value = '1.299'
$ws.Range("E30").Value = '  -10.90%  '
$ws.Range("D31").NumberFormat = '@'
$ws.Range("D31").Value = '132.86'
$ws.Range("E31").Value = '  -1.23%  '
$ws.Range("D32").NumberFormat = '@'
$ws.Range("D32").Value = '4.065'
$ws.Range("E32").Value = '  +0.77%  '
$ws.Range("D33").NumberFormat = '@'
$ws.Range("D33").Value = '5.960'
$ws.Range("E33").Value = '  -1.36%  '
$ws.Range("D34").NumberFormat = '@'
$ws.Range("D34").Value = '0.09231'
$ws.Range("E34").Value = '  +3.71%  '
$ws.Range("D35").NumberFormat = '@'
$ws.Range("D35").Value = '12.37'
$ws.Range("E35").Value = '  -2.44%  '
$ws.Range("D36").NumberFormat = '@'
$ws.Range("D36").Value = '0.02355'
$ws.Range("E36").Value = '  -2.09%  '
$ws.Range("D37").NumberFormat = '@'
$ws.Range("D37").Value = '0.6718'
$ws.Range("E37").Value = '  -1.44%  '
$ws.Range("D38").NumberFormat = '@'
$ws.Range("D38").Value = '5.245'
$ws.Range("D39").NumberFormat = '@'
$ws.Range("D39").Value = '0.06270'
$ws.Range("E39").Value = '  -1.33%  '
$ws.Range("D40").NumberFormat = '@'
$ws.Range("D40").Value = '0.2175'
$ws.Range("E40").Value = '  +1.01%  '
$ws.Range("E41").Value = '  -0.09%  '
$ws.Range("D42").NumberFormat = '@'
$ws.Range("D42").Value = '1.216'
$ws.Range("E42").Value = '  -1.62%  '
$ws.Range("D43").NumberFormat = '@'
$ws.Range("D43").Value = '8.086'
$ws.Range("E43").Value = '  -1.12%  '
$ws.Range("D44").NumberFormat = '@'
$ws.Range("D44").Value = '0.9998'
$ws.Range("E44").Value = '  -0.10%  '
$ws.Range("D45").NumberFormat = '@'
$ws.Range("D45").Value = '14.07'
$ws.Range("E45").Value = '  -0.65%  '
$ws.Range("E46").Value = '  +0.21%  '
$ws.Range("D47").NumberFormat = '@'
$ws.Range("D47").Value = '0.6126'
$ws.Range("E47").Value = '  -2.02%  '
$ws.Range("D48").NumberFormat = '@'
$ws.Range("D48").Value = '128.16'
$ws.Range("E48").Value = '  -3.35%  '
$ws.Range("D49").NumberFormat = '@'
$ws.Range("D49").Value = '2.048'
$ws.Range("E49").Value = '  -0.88%  '
$ws.Range("D50").NumberFormat = '@'
$ws.Range("D50").Value = '1.180'
$ws.Range("E50").Value = '  -2.14%  '
$ws.Range("D51").NumberFormat = '@'
$ws.Range("D51").Value = '0.07097'
$ws.Range("E51").Value = '  -5.30%  '
